# Peer review in 3 columns. Fixed sudden overlapping of 'with' entries.
# Inserts a new grant/funding record (Milena Vasquez-Amezquita project) at
# the top of the "education" sheet's grant table, pushing the existing
# four records down by three rows each.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Insert 3 blank rows above the current row 2 (which holds the first
# existing grant record) to make room for the new entry + its two
# "with"-style sub-rows (PI line + amount line).
$ws.Rows("2:4").Insert()

# --- New grant record (row 2: what / when / with / where / why) ---
$ws.Range("A2").Value = "XI \href{https://www.unbosque.edu.co/investigaciones/convocatorias-investigacion}{Internal Call for Financing Research and Technological Innovation Projects El Bosque University}, 2023"
$ws.Range("B2").Value = "Feb. 2024 - Present"
$ws.Range("C2").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}"
$ws.Range("D2").Value = "Bogota, Colombia"
$ws.Range("E2").Value = "Project: \textit{Effect of real and simulated resource control on androphilic women's preferences for masculinity in men's faces: an experimental study using eye-tracking}"
$ws.Rows(2).RowHeight = 60

# --- Sub-row 3: Principal Investigator line ---
$ws.Range("E3").Value = "PI: \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}"

# --- Sub-row 4: grant amount line (currency-style number format) ---
$ws.Range("E4").Value = "COP\$89.979.750"
$ws.Range("E4").Style = $ws.Range("E6").Style

# Move the active selection to B4, matching where the author's cursor
# ended up after entering the new record.
$ws.Range("B4").Select()
